$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 and IF headers in column I and J, matching style of existing header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-12 for columns I (I0) and J (IF)
$data = @{
    2  = @(1, 3)
    3  = @(1, 4)
    4  = @(1, 3)
    5  = @(1, 5)
    6  = @(1, 7)
    7  = @(1, 5)
    8  = @(1, 5)
    9  = @(9, 9)
    10 = @(8, 8)
    11 = @(9, 9)
    12 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
